# MIC supply filter fix.
# Re-balances which capacitors/resistors fall into which BOM "value" bucket,
# renames the ferrite bead part, re-numbers a couple of capacitor designators
# (a new C27 appears), adds a price-column currency format + a total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: ferrite bead value changed from 1k to 600 ---
$ws.Range("B12").Value = "600 ferrite bead"

# --- Row 14: 1n-50V C0G 1% capacitor designator list grows (C9 removed, C14 added) ---
$ws.Range("D14").Value = "C1, C2, C3, C4, C5, C6, C10, C11, C12, C13, C14"

# --- Rows 15-18: capacitor value rows re-ordered & designators re-balanced ---
$ws.Range("A15").Formula = "=`$M`$1*2"
$ws.Range("B15").Value = "1u-16V"
$ws.Range("D15").Value = "C7, C23"
$ws.Range("E15").Value = "X7R"

$ws.Range("A16").Formula = "=`$M`$1*8"
$ws.Range("B16").Value = "100n-50V"
$ws.Range("D16").Value = "C8, C9, C15, C16, C17, C20, C22, C26"
$ws.Range("E16").Value = "X7R"

$ws.Range("A17").Formula = "=`$M`$1*2"
$ws.Range("B17").Value = "22p-50V"
$ws.Range("D17").Value = "C18, C19"
$ws.Range("E17").Value = "C0G"

$ws.Range("A18").Formula = "=`$M`$1*3"
$ws.Range("B18").Value = "10u-16V"
$ws.Range("D18").Value = "C21, C25, C27"
$ws.Range("E18").Value = "X5R"

# --- Row 19: 10n-50V capacitor designator updated (C23 -> C24) ---
$ws.Range("D19").Value = "C24"

# --- Apply a Forint currency number format to the whole PRICE column ---
$ws.Range("I2:I44").NumberFormat = "#,##0 ""Ft"""

# --- Add a total row under the table that sums the price column ---
$ws.Range("I2").Copy()
$ws.Range("I45").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("I45").Formula = "=SUM(I2:I44)"

# --- Misc view bookkeeping to mirror the authored workbook state ---
# (PART column got a bit wider once the designator lists grew)
$ws.Columns.Item(4).ColumnWidth = 36.6

$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("I40").Select()
